# Fix an alignment problem with the report
#
# 1) Remove the stray empty paragraph that sits right before the
#    "Green State:" bulleted list item (directly after "...the following
#    key states:").
# 2) The "Red State:" paragraph text had been split across two runs at a
#    page-break boundary ("...pedestrians " | "and vehicles..."). Merge
#    them back into a single, contiguous sentence.
# 3) The <w:lastRenderedPageBreak/> marker that used to sit on the
#    "and vehicles..." run now belongs at the start of the following
#    "Manual Override:" run (the page now breaks there instead).

$d = $word.ActiveDocument

# --- Step 1: delete the empty paragraph before the "Green State:" item ---
$searchRange = $d.Content
$searchRange.Find.Execute("the following key states:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingEnd = $searchRange.End

$greenRange = $d.Content
$greenRange.Find.Execute("Green State:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$greenStart = $greenRange.Start

$emptyPara = $d.Range($headingEnd, $greenStart)
$emptyPara.Delete()

# --- Step 2: merge the split "pedestrians " / "and vehicles..." runs ---
$d.Content.Find.Execute("pedestrians and vehicles on the secondary road", $true, $false, $false, $false, $false, $true, 1, $false, "pedestrians and vehicles on the secondary road", 2) | Out-Null

# --- Step 3: move <w:lastRenderedPageBreak/> onto the "Manual Override:" run ---
$moRange = $d.Content
$moRange.Find.Execute("Manual Override:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$moStart = $moRange.Start

$insPt = $d.Range($moStart, $moStart)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="ar-EG"/></w:rPr><w:lastRenderedPageBreak/><w:t>Manual Override:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPt.InsertXML($xml)

$dupRange = $d.Range($moStart + 16, $moStart + 32)
$dupRange.Delete()
